$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(8,8),
    @(6,7),
    @(9,9),
    @(9,9),
    @(7,7),
    @(9,9),
    @(5,6),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,6),
    @(9,9),
    @(7,8),
    @(6,6),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8)
)

for ($i = 0; $i -lt 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
